# Edit summary (per the supplied OOXML diff):
#   1. The table on slide 16 (the "PLENARY - COMPLETE THE MISSING GAPS" slide)
#      gets its table style switched from the custom "{89929DD1-...}" style to
#      the built-in "{EF840ED4-BAAF-45B2-88D6-8397BBE94F02}" (Medium Style 2 -
#      Accent 1) table style.
#   2. The deck's theme colour scheme (the one backing slideMaster1 / every
#      slide, stored in ppt/theme/theme2.xml) is switched from the custom
#      "Integral" palette back to the standard Office palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{EF840ED4-BAAF-45B2-88D6-8397BBE94F02}")

# --- 2) Theme colours: Integral -> Office ----------------------------------
# Order of ThemeColorScheme.Colors() is:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# Values below are the standard "Office" theme RGB values expressed the way
# VBA's RGB() macro would encode them (R + G*256 + B*65536).
$officeColors = @(
    0,           # dk1     000000
    16777215,    # lt1     FFFFFF
    6968388,     # dk2     44546A
    15132391,    # lt2     E7E6E6
    13998939,    # accent1 5B9BD5
    3243501,     # accent2 ED7D31
    10855845,    # accent3 A5A5A5
    49407,       # accent4 FFC000
    12874308,    # accent5 4472C4
    4697456,     # accent6 70AD47
    12673797,    # hlink   0563C1
    7491477      # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
